$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Películas")
$tbl = $ws.ListObjects.Item(1)
$newRow = $tbl.ListRows.Add()
$rowRange = $newRow.Range

$srcRange = $ws.Range("B120:I120")
$srcRange.Copy()
$rowRange.PasteSpecial(-4122)

$rowRange.Item(1).Value = "Troll 2"
$rowRange.Item(2).Formula = "=AVERAGE(D121,E121,E121,F121,G121,H121,H121,I121)"
$rowRange.Item(3).Value = 7
$rowRange.Item(4).Value = 7
$rowRange.Item(5).Value = 7
$rowRange.Item(6).Value = 7
$rowRange.Item(7).Value = 5.5
$rowRange.Item(8).Value = 4.9

$rng = $tbl.Range
$rng.Sort($ws.Range("C3"), 2, $ws.Range("E3"), , 1, , , 1)

# Now find the B-cell that has the "Ladrones con clase" text (formerly styled s=9) and force its format
# to match a typical B-cell format (e.g. copy format from B80, which uses s=2).
$fixSrc = $ws.Range("B80")
$fixSrc.Copy()
$ws.Range("B117").PasteSpecial(-4122)
Write-Output "done"
